$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AE1 - copy formatting from AD1 (same style as N1:AD1), then set the date text
$ws.Range("AD1").Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = "16-10-2020"

# Data cells AE2:AE36
$ws.Range("AE2").Value = 190
$ws.Range("AE3").Value = 40047
$ws.Range("AE4").Value = 3052
$ws.Range("AE5").Value = 28804
$ws.Range("AE6").Value = 11038
$ws.Range("AE7").Value = 1044
$ws.Range("AE8").Value = 28187
$ws.Range("AE9").Value = 71
$ws.Range("AE10").Value = 22605
$ws.Range("AE11").Value = 4084
$ws.Range("AE12").Value = 14782
$ws.Range("AE13").Value = 10364
$ws.Range("AE14").Value = 2654
$ws.Range("AE15").Value = 9058
$ws.Range("AE16").Value = 6892
$ws.Range("AE17").Value = 113557
$ws.Range("AE18").Value = 94609
$ws.Range("AE19").Value = 1018
$ws.Range("AE20").Value = 14157
$ws.Range("AE21").Value = 192936
$ws.Range("AE22").Value = 3193
$ws.Range("AE23").Value = 2445
$ws.Range("AE24").Value = 108
$ws.Range("AE25").Value = 1453
$ws.Range("AE26").Value = 22387
$ws.Range("AE27").Value = 4551
$ws.Range("AE28").Value = 7090
$ws.Range("AE29").Value = 21587
$ws.Range("AE30").Value = 312
$ws.Range("AE31").Value = 41872
$ws.Range("AE32").Value = 23315
$ws.Range("AE33").Value = 3105
$ws.Range("AE34").Value = 5682
$ws.Range("AE35").Value = 36295
$ws.Range("AE36").Value = 31984
